# Applies the "Added more to developer sheet and created tables" commit.
#
# Summary of the change (per the OOXML diff):
#  1. developers sheet (rows 61-100): fill in the Location (C) and Founded (D)
#     columns that were previously blank for those developers.
#       - C81 picks up a "d-mmm-yy" number format (same numFmtId as the
#         pre-existing style index 3).
#       - C88 / C95 pick up word-wrap (a brand-new style, numFmtId 0 +
#         wrapText, which becomes cellXfs index 5).
#       - D90 holds a literal text date ("23 September 1889") rather than a
#         serial date number.
#  2. Selection / active-cell bookkeeping changed on every sheet, and the
#     "developers" sheet becomes the active tab (activeTab goes from the
#     last sheet, "publishers", to "developers").

$wb = $excel.ActiveWorkbook

$wsDevelopers = $wb.Worksheets.Item("developers")

# ---------------------------------------------------------------------
# 1. Populate the Location (C) / Founded (D) columns for rows 61-100.
# ---------------------------------------------------------------------
$ws = $wsDevelopers

$ws.Range("C61").Value = "Milan, Italy"
$ws.Range("D61").Value = 1998

$ws.Range("C62").Value = "Montreal, Canada"
$ws.Range("D62").Value = 39412

$ws.Range("C63").Value = "Toronto, Canada"
$ws.Range("D63").Value = 40422

$ws.Range("C64").Value = "Montreuil, France"
$ws.Range("D64").Value = 31499

$ws.Range("C65").Value = "Frankfurt, Germany"
$ws.Range("D65").Value = 36404

$ws.Range("C66").Value = "Brighton, England"
$ws.Range("D66").Value = 2004

$ws.Range("C67").Value = "Dallas, Texas, U.S."
$ws.Range("D67").Value = 1994

$ws.Range("C68").Value = "Madison, Wisconsin, U.S."
$ws.Range("D68").Value = 35704

$ws.Range("C69").Value = "Japan"
$ws.Range("D69").Value = 36708

$ws.Range("C70").Value = "Kirkland, Washington, United States"
$ws.Range("D70").Value = 34632

$ws.Range("C71").Value = "Bellevue, Washington, US"
$ws.Range("D71").Value = 35301

$ws.Range("C72").Value = "Stockholm, Sweden"
$ws.Range("D72").Value = 1998

$ws.Range("C73").Value = "Carlsbad, California, US"
$ws.Range("D73").Value = 30682

$ws.Range("C74").Value = "Tokyo, Japan"
$ws.Range("D74").Value = 30590

$ws.Range("C75").Value = "Redmond, Washington, US"
$ws.Range("D75").Value = 2007

$ws.Range("C76").Value = "Ginza, Chūō, Tokyo, Japan"
$ws.Range("D76").Value = 25283

$ws.Range("C77").Value = "Chicago, US"
$ws.Range("D77").Value = 40288

$ws.Range("C78").Value = "Redmond, Washington, US"
$ws.Range("D78").Value = 2009

$ws.Range("C79").Value = "Liverpool, England"
$ws.Range("D79").Value = 1988

# Row 80 only gets a Founded date, no Location.
$ws.Range("D80").Value = 34095

$ws.Range("C81").Value = "Yodogawa-ku, Osaka, Japan"
$ws.Range("C81").NumberFormat = "d-mmm-yy"
$ws.Range("D81").Value = 33312

$ws.Range("C82").Value = "Japan"
$ws.Range("D82").Value = 1988

$ws.Range("C83").Value = "Emeryville, California, US"
$ws.Range("D83").Value = 1992

$ws.Range("C84").Value = "Greater New York Area, East Coast, Northeastern US"
$ws.Range("D84").Value = 1995

$ws.Range("C85").Value = "Shanghai, China"
$ws.Range("D85").Value = 36678

$ws.Range("C86").Value = "Lemont, Illinois"
$ws.Range("D86").Value = 1993

$ws.Range("C87").Value = "Manhattan Beach, California, US"
$ws.Range("D87").Value = 36403

$ws.Range("C88").Value = "Los Angeles, California"
$ws.Range("C88").WrapText = $true
$ws.Range("D88").Value = 1997

$ws.Range("C89").Value = "San Mateo, California"
$ws.Range("D89").Value = 1994

$ws.Range("D90").Value = "23 September 1889"
$ws.Range("C90").Value = "11-1 Kamitoba Hokodatecho, Minami-ku, Kyoto, Japan"

$ws.Range("C91").Value = "Sheffield, England"
$ws.Range("D91").Value = 1977

$ws.Range("C92").Value = "Tokyo, Japan"
$ws.Range("D92").Value = 38777

$ws.Range("C93").Value = "Shinagawa City, Tokyo, Japan"
$ws.Range("D93").Value = 35156

$ws.Range("C94").Value = "Wellington, New Zealand"
$ws.Range("D94").Value = 35551

$ws.Range("C95").Value = "Kyoto, Japan"
$ws.Range("C95").WrapText = $true
$ws.Range("D95").Value = 2004

$ws.Range("C96").Value = "Shalford, United Kingdom"
$ws.Range("D96").Value = 1998

$ws.Range("C97").Value = "Austin, Texas, U.S."
$ws.Range("D97").Value = 36059

$ws.Range("C98").Value = "New Taipei, Taiwan"
$ws.Range("D98").Value = 32813

$ws.Range("C99").Value = "Suginami, Tokyo, Japan"
$ws.Range("D99").Value = 27870

$ws.Range("C100").Value = "Ageo, Saitama, Japan"
$ws.Range("D100").Value = 29403

# ---------------------------------------------------------------------
# 2. Window / selection bookkeeping on every sheet. The last sheet that
#    is Activate()-d ends up as the active tab, so "developers" (which
#    becomes the active tab in the target workbook) is activated last.
# ---------------------------------------------------------------------
$wsVideoGames = $wb.Worksheets.Item("video_games")
$wsVideoGames.Activate()
$wsVideoGames.Range("G19").Select()

$wsGenres = $wb.Worksheets.Item("genres")
$wsGenres.Activate()
$wsGenres.Range("A12").Select()

$wsPublishers = $wb.Worksheets.Item("publishers")
$wsPublishers.Activate()
$wsPublishers.Range("D14").Select()

$wsDevelopers.Activate()
$wsDevelopers.Range("D101").Select()
